# Auto-generated edit script applying numeric updates described by the diff.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")  # row 17
$ws.Range("H17").Value = 968.65515
$ws.Range("J17").Value = 1095.2609
$ws.Range("L17").Value = 3285.7827
$ws.Range("N17").Value = -3621.7827

$ws = $wb.Worksheets.Item("ALC")  # row 46
$ws.Range("H46").Value = 20833922
$ws.Range("J46").Value = 733.3333
$ws.Range("L46").Value = 2199.9999
$ws.Range("N46").Value = -2437.9999

$ws = $wb.Worksheets.Item("ALC")  # row 58
$ws.Range("H58").Value = 2099.6072
$ws.Range("I58").Value = 366.27274
$ws.Range("J58").Value = 3221.1765
$ws.Range("K58").Value = 1098.81822
$ws.Range("L58").Value = 9663.529500000001
$ws.Range("M58").Value = -948.8182200000001
$ws.Range("N58").Value = -9963.529500000001

$ws = $wb.Worksheets.Item("ALC")  # row 59
$ws.Range("H59").Value = 750
$ws.Range("J59").Value = 0
$ws.Range("L59").Value = 0
$ws.Range("N59").ClearContents()

$ws = $wb.Worksheets.Item("ALC")  # row 60
$ws.Range("H60").Value = 20833922
$ws.Range("J60").Value = 733.3333
$ws.Range("L60").Value = 2199.9999
$ws.Range("N60").Value = -3167.9999

$ws = $wb.Worksheets.Item("ARM")  # row 2
$ws.Range("H2").Value = 1285.1818
$ws.Range("I2").Value = 1527.75
$ws.Range("J2").Value = 1146.5714
$ws.Range("K2").Value = 1527.75
$ws.Range("L2").Value = 1146.5714
$ws.Range("M2").Value = -1414.75
$ws.Range("N2").Value = -1372.5714

$ws = $wb.Worksheets.Item("ARM")  # row 32
$ws.Range("H32").Value = 5638.99
$ws.Range("I32").Value = 5067.358
$ws.Range("J32").Value = 16500
$ws.Range("K32").Value = 5067.358
$ws.Range("L32").Value = 16500
$ws.Range("M32").Value = -4780.358
$ws.Range("N32").Value = -17074

$ws = $wb.Worksheets.Item("ARM")  # row 74
$ws.Range("H74").Value = 1448.4517
$ws.Range("I74").Value = 1299.9
$ws.Range("K74").Value = 1299.9
$ws.Range("M74").Value = -425.9000000000001

$ws = $wb.Worksheets.Item("ARM")  # row 77
$ws.Range("H77").Value = 1448.4517
$ws.Range("I77").Value = 1299.9
$ws.Range("K77").Value = 6499.5
$ws.Range("M77").Value = -2131.5

$ws = $wb.Worksheets.Item("ARM")  # row 116
$ws.Range("H116").Value = 1285.1818
$ws.Range("I116").Value = 1527.75
$ws.Range("J116").Value = 1146.5714
$ws.Range("K116").Value = 1527.75
$ws.Range("L116").Value = 1146.5714
$ws.Range("M116").Value = 766.25
$ws.Range("N116").Value = -5734.5714

$ws = $wb.Worksheets.Item("ARM")  # row 122
$ws.Range("H122").Value = 1119.5
$ws.Range("I122").Value = 1023.3333
$ws.Range("J122").Value = 1696.5
$ws.Range("K122").Value = 3069.9999
$ws.Range("L122").Value = 5089.5
$ws.Range("M122").Value = -619.9998999999998
$ws.Range("N122").Value = -9989.5

$ws = $wb.Worksheets.Item("ARM")  # row 132
$ws.Range("H132").Value = 1671284.4
$ws.Range("I132").Value = 4198.442
$ws.Range("J132").Value = 3911431.2
$ws.Range("K132").Value = 12595.326
$ws.Range("L132").Value = 11734293.6
$ws.Range("M132").Value = -10065.326
$ws.Range("N132").Value = -11739353.6

$ws = $wb.Worksheets.Item("BSM")  # row 3
$ws.Range("H3").Value = 1285.1818
$ws.Range("I3").Value = 1527.75
$ws.Range("J3").Value = 1146.5714
$ws.Range("K3").Value = 1527.75
$ws.Range("L3").Value = 1146.5714
$ws.Range("M3").Value = -1413.75
$ws.Range("N3").Value = -1374.5714

$ws = $wb.Worksheets.Item("BSM")  # row 43
$ws.Range("H43").Value = 231020
$ws.Range("J43").Value = 231020
$ws.Range("L43").Value = 231020
$ws.Range("N43").Value = -231382

$ws = $wb.Worksheets.Item("BSM")  # row 134
$ws.Range("H134").Value = 2136.597
$ws.Range("I134").Value = 1325.5264
$ws.Range("J134").Value = 3199.3794
$ws.Range("K134").Value = 3976.5792
$ws.Range("L134").Value = 9598.138199999999
$ws.Range("M134").Value = -1441.5792
$ws.Range("N134").Value = -14668.1382

$ws = $wb.Worksheets.Item("CRP")  # row 19
$ws.Range("H19").Value = 431
$ws.Range("I19").Value = 431
$ws.Range("J19").Value = 0
$ws.Range("K19").Value = 431
$ws.Range("L19").Value = 0
$ws.Range("M19").ClearContents()
$ws.Range("N19").Value = -261

$ws = $wb.Worksheets.Item("CRP")  # row 24
$ws.Range("H24").Value = 431
$ws.Range("I24").Value = 431
$ws.Range("J24").Value = 0
$ws.Range("K24").Value = 431
$ws.Range("L24").Value = 0
$ws.Range("M24").ClearContents()
$ws.Range("N24").Value = -261

$ws = $wb.Worksheets.Item("CRP")  # row 58
$ws.Range("H58").Value = 3147.92
$ws.Range("I58").Value = 2011.8
$ws.Range("J58").Value = 4852.1
$ws.Range("K58").Value = 2011.8
$ws.Range("L58").Value = 4852.1
$ws.Range("M58").Value = -1808.8
$ws.Range("N58").Value = -5258.1

$ws = $wb.Worksheets.Item("CRP")  # row 106
$ws.Range("H106").Value = 35000.5
$ws.Range("J106").Value = 35000.5
$ws.Range("L106").Value = 35000.5
$ws.Range("N106").Value = -37524.5

$ws = $wb.Worksheets.Item("CRP")  # row 122
$ws.Range("H122").Value = 111112296
$ws.Range("I122").Value = 142857810
$ws.Range("J122").Value = 3007
$ws.Range("K122").Value = 428573430
$ws.Range("L122").Value = 9021
$ws.Range("M122").Value = -428570980
$ws.Range("N122").Value = -13921

$ws = $wb.Worksheets.Item("CRP")  # row 134
$ws.Range("H134").Value = 1664.4
$ws.Range("I134").Value = 1161.0588
$ws.Range("J134").Value = 2322.6155
$ws.Range("K134").Value = 3483.1764
$ws.Range("L134").Value = 6967.8465
$ws.Range("M134").Value = -948.1764000000003
$ws.Range("N134").Value = -12037.8465

$ws = $wb.Worksheets.Item("CRP")  # row 136
$ws.Range("H136").Value = 3147.92
$ws.Range("I136").Value = 2011.8
$ws.Range("J136").Value = 4852.1
$ws.Range("K136").Value = 6035.4
$ws.Range("L136").Value = 14556.3
$ws.Range("M136").Value = -3485.4
$ws.Range("N136").Value = -19656.3

$ws = $wb.Worksheets.Item("CUL")  # row 44
$ws.Range("H44").Value = 350.66666
$ws.Range("I44").Value = 323.1111
$ws.Range("J44").Value = 433.33334
$ws.Range("K44").Value = 969.3333
$ws.Range("L44").Value = 1300.00002
$ws.Range("M44").Value = -571.3333
$ws.Range("N44").Value = -2096.00002

$ws = $wb.Worksheets.Item("CUL")  # row 46
$ws.Range("H46").Value = 586.75
$ws.Range("I46").Value = 579.125
$ws.Range("J46").Value = 602
$ws.Range("K46").Value = 1737.375
$ws.Range("L46").Value = 1806
$ws.Range("M46").Value = -1646.375
$ws.Range("N46").Value = -1988

$ws = $wb.Worksheets.Item("CUL")  # row 59
$ws.Range("H59").Value = 1809.8
$ws.Range("I59").Value = 775
$ws.Range("J59").Value = 2499.6667
$ws.Range("K59").Value = 2325
$ws.Range("L59").Value = 7499.000100000001
$ws.Range("M59").Value = -1785
$ws.Range("N59").Value = -8579.000100000001

$ws = $wb.Worksheets.Item("CUL")  # row 61
$ws.Range("H61").Value = 442.27274
$ws.Range("I61").Value = 98
$ws.Range("J61").Value = 855.4
$ws.Range("K61").Value = 294
$ws.Range("L61").Value = 2566.2
$ws.Range("M61").Value = -79
$ws.Range("N61").Value = -2996.2

$ws = $wb.Worksheets.Item("GSM")  # row 9
$ws.Range("H9").Value = 616.875
$ws.Range("I9").Value = 616.875
$ws.Range("K9").Value = 616.875
$ws.Range("M9").Value = -446.875

$ws = $wb.Worksheets.Item("GSM")  # row 10
$ws.Range("H10").Value = 300
$ws.Range("I10").Value = 300
$ws.Range("K10").Value = 300
$ws.Range("M10").Value = -131

$ws = $wb.Worksheets.Item("GSM")  # row 107
$ws.Range("H107").Value = 949.8333
$ws.Range("I107").Value = 945.2727
$ws.Range("J107").Value = 1000
$ws.Range("K107").Value = 945.2727
$ws.Range("L107").Value = 1000
$ws.Range("M107").Value = 974.7273
$ws.Range("N107").Value = -4840

$ws = $wb.Worksheets.Item("GSM")  # row 126
$ws.Range("H126").Value = 3908168.5
$ws.Range("I126").Value = 6251574.5
$ws.Range("K126").Value = 18754723.5
$ws.Range("M126").Value = -18752253.5

$ws = $wb.Worksheets.Item("GSM")  # row 132
$ws.Range("H132").Value = 2363.83
$ws.Range("I132").Value = 1368.4667
$ws.Range("J132").Value = 3662.1304
$ws.Range("K132").Value = 4105.4001
$ws.Range("L132").Value = 10986.3912
$ws.Range("M132").Value = -1575.4001
$ws.Range("N132").Value = -16046.3912

$ws = $wb.Worksheets.Item("LTW")  # row 40
$ws.Range("H40").Value = 20835926
$ws.Range("I40").Value = 33335522
$ws.Range("J40").Value = 3267.389
$ws.Range("K40").Value = 33335522
$ws.Range("L40").Value = 3267.389
$ws.Range("M40").Value = -33335386
$ws.Range("N40").Value = -3539.389

$ws = $wb.Worksheets.Item("LTW")  # row 105
$ws.Range("H105").Value = 28900
$ws.Range("J105").Value = 28900
$ws.Range("L105").Value = 28900
$ws.Range("N105").Value = -35888

$ws = $wb.Worksheets.Item("LTW")  # row 132
$ws.Range("H132").Value = 22698.818
$ws.Range("I132").Value = 30669.053
$ws.Range("J132").Value = 4883
$ws.Range("K132").Value = 92007.159
$ws.Range("L132").Value = 14649
$ws.Range("M132").Value = -89477.159
$ws.Range("N132").Value = -19709

$ws = $wb.Worksheets.Item("WVR")  # row 108
$ws.Range("H108").Value = 40000
$ws.Range("J108").Value = 40000
$ws.Range("L108").Value = 40000
$ws.Range("N108").Value = -47680
